# edit.ps1 -- applies the "Quantum Entanglement" -> "History" rewrite
# described by the commit diff, via Word COM-interop calls.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Find-Range($text) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng
}

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
Replace-Text "Quantum Entanglement's Secrets Unveiled" "A Journey Through Time and Space: The Captivating Story of History"

# ---------------------------------------------------------------------
# Byline: "Jonas Hoffstadter" -> "Dr. Albert Meadows"
# ---------------------------------------------------------------------
Replace-Text "Jonas Hoffstadter" "Dr. Albert Meadows"

# ---------------------------------------------------------------------
# Email address
# ---------------------------------------------------------------------
Replace-Text "jonashoffstadter@protonmail" "albertmeadows@academy"
$emailPara = $d.Paragraphs(3).Range
$emailPara.Find.Execute("com", $true, $true, $false, $false, $false, $true, 1, $false, "org", 2) | Out-Null

# ---------------------------------------------------------------------
# Intro paragraph (paragraph 5)
# ---------------------------------------------------------------------
Replace-Text "In the realm of quantum mechanics lies a phenomenon that has ignited curiosity and perplexity in equal measure - quantum entanglement" "Our world is a tapestry of stories whispered through time, each era weaving its unique thread into the grand narrative of history"

Replace-Text " As particles pair and share an indefinable link, scientists strive to unravel the mysteries of this elusive connection, venturing into uncharted territories of physics" " As students of history, we embark on a grand adventure, unveiling the enigmatic tapestry of past events and unraveling the enigmas of civilizations long gone"

Replace-Text " In this exploration, we will venture into the intricate depths of quantum entanglement, demystifying its paradoxical behaviors and unveiling its fundamental implications for our understanding of reality" " We uncover the hidden narratives concealed within ancient ruins, deciphering the secrets etched in stone and parchment"

$r = Find-Range " We uncover the hidden narratives concealed within ancient ruins, deciphering the secrets etched in stone and parchment"
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" Through the symphony of history, we dance with forgotten heroes and heroines, their triumphs and struggles resonating with our present lives")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
# the run that used to be the lone "." right after the original run4 text
# becomes the next sentence; it is the single character immediately
# following the text we just inserted.
$period = $d.Range($r.End, $r.End + 1)
$period.Text = " History is not merely a collection of dates and names; it is a vibrant chronicle of human experience, a kaleidoscope of cultures, and a symphony of civilizations"

Replace-Text "From the early postulates of Schrodinger's cat to the theoretical framework of Bell's inequality, we will delve into the history and scientific milestones that have shaped our comprehension of entanglement" "We traverse through the ages, witnessing the rise and fall of empires, exploring the evolution of ideas, and encountering the architects of destiny"

Replace-Text " Through the lens of experimental evidence, we will witness the remarkable feats achieved in laboratories across the globe, demonstrating the baffling correlations between entangled particles that appear to transcend the limitations of time and space" " In the annals of history, we find solace in the wisdom of philosophers, navigate the treacherous waters of political intrigue, and unravel the mysteries of scientific discoveries"

$r2 = Find-Range " In the annals of history, we find solace in the wisdom of philosophers, navigate the treacherous waters of political intrigue, and unravel the mysteries of scientific discoveries"
$r2.Collapse(0)
$r2.InsertAfter(".")
$r2.Collapse(0)
$r2.InsertAfter(" It is a journey through time and space, connecting us to the past and propelling us into the future")
$r2.Collapse(0)
$r2.InsertAfter(".")
$r2.Collapse(0)
$period2 = $d.Range($r2.End, $r2.End + 1)
$period2.Text = " History is a mirror that reflects our collective heritage, a roadmap guiding our present, and a window into the possibilities of tomorrow"

Replace-Text "Moreover, we will delve into the profound implications of entanglement, from its potential applications in quantum computation to its challenge to our classical notions of locality and determinism" "History provides invaluable lessons for the future, teaching us the ramifications of human choices, the cyclical nature of conflict and cooperation, and the enduring power of human resilience"

Replace-Text " As we navigate the theoretical and practical frontiers of quantum entanglement, we will ponder the very nature of reality, questioning the fundamental building blocks of the universe and the interconnectedness of all things" " It challenges us to confront our triumphs and failures, forging a path forward that honors the sacrifices of those who came before us"

$r3 = Find-Range " It challenges us to confront our triumphs and failures, forging a path forward that honors the sacrifices of those who came before us"
$r3.Collapse(0)
$r3.InsertAfter(".")
$r3.Collapse(0)
$r3.InsertAfter(" As we delve deeper into the story of humanity, we cultivate empathy, compassion, and a profound sense of interconnectedness, recognizing that our shared past binds us together in the tapestry of time")

# ---------------------------------------------------------------------
# Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------
Replace-Text "Our journey into the enigmatic world of quantum entanglement has illuminated the remarkable phenomena and profound implications of this perplexing linkage between particles" "History is a captivating narrative that unveils the enigmas of the past and shapes our understanding of the present"

Replace-Text " From the theoretical constructs of Schrodinger and Bell to the experimental breakthroughs in laboratories, we have witnessed the mesmerizing behaviors of entangled particles and their defiance of classical physics" " By exploring the annals of history, we traverse through time and "

$r4 = Find-Range " By exploring the annals of history, we traverse through time and "
$r4.Collapse(0)
$r4.InsertAfter("space, witnessing the rise and fall of empires, deciphering forgotten wisdom, and uncovering the secrets of civilizations long gone")

Replace-Text " The potential for quantum entanglement to revolutionize fields such as communication, computation, and cryptography holds immense promise" " History provides invaluable lessons, teaching us the consequences of human choices and the power of resilience"

Replace-Text " Yet, it also invites us to confront the unsettling questions about the nature of reality and our place within it" " As we delve deeper into the tapestry of history, we cultivate empathy, compassion, and a sense of interconnectedness, recognizing that our shared past binds us together"

Replace-Text " As we continue to unravel the secrets of quantum entanglement, we venture ever closer to the enigmatic core of existence, where the boundaries of our knowledge blend with the infinite mysteries of the universe" " History is not just a collection of dates and names; it is a vibrant symphony of human experience that resonates with our lives and propels us into the future"

# ---------------------------------------------------------------------
# New trailing empty paragraph
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

Write-Output "edit complete"
